# Update the underlying data tables that back the "Bar chart" and "Pi chart"
# charts. Both charts plot directly from these worksheet ranges, so editing
# the cells here is the correct, Excel-faithful way to change what the
# charts show.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# "Bar chart" sheet - table in A1:E5 (columns B:E hold TP/FP/TN/FN counts
# for each of the four rows: k-NN Regression, k-NN Classification,
# Neural Network, Linear Regression).
# ---------------------------------------------------------------------
$wsBar = $wb.Worksheets.Item("Bar chart")

$wsBar.Range("B2").Value = 8028
$wsBar.Range("C2").Value = 7998
$wsBar.Range("D2").Value = 5515
$wsBar.Range("E2").Value = 3437

$wsBar.Range("B3").Value = 64
$wsBar.Range("C3").Value = 71
$wsBar.Range("D3").Value = 2303
$wsBar.Range("E3").Value = 0

$wsBar.Range("B4").Value = 3537
$wsBar.Range("C4").Value = 3530
$wsBar.Range("D4").Value = 1298
$wsBar.Range("E4").Value = 3601

$wsBar.Range("B5").Value = 41
$wsBar.Range("C5").Value = 61
$wsBar.Range("D5").Value = 961
$wsBar.Range("E5").Value = 0

# ---------------------------------------------------------------------
# "Pi chart" sheet - table in A1:E5 (columns B:E hold the TP/FP/TN/FN
# values for each of the four rows: k-NN Regression, k-NN Classification,
# Neural Network, Linear Regression).
# ---------------------------------------------------------------------
$wsPie = $wb.Worksheets.Item("Pi chart")

$wsPie.Range("B2").Value = 8028
$wsPie.Range("C2").Value = 64
$wsPie.Range("D2").Value = 3537
$wsPie.Range("E2").Value = 41

$wsPie.Range("B3").Value = 7998
$wsPie.Range("C3").Value = 71
$wsPie.Range("D3").Value = 3530
$wsPie.Range("E3").Value = 61

$wsPie.Range("B4").Value = 5515
$wsPie.Range("C4").Value = 2303
$wsPie.Range("D4").Value = 1298
$wsPie.Range("E4").Value = 961

$wsPie.Range("B5").Value = 3437
$wsPie.Range("C5").Value = 0
$wsPie.Range("D5").Value = 3601
$wsPie.Range("E5").Value = 0

# ---------------------------------------------------------------------
# Update the selected cell on each sheet to match the saved view state,
# then leave the "Bar chart" sheet active/selected (as in the workbook).
# ---------------------------------------------------------------------
$wsPie.Activate()
$wsPie.Range("F11").Select()

$wsBar.Activate()
$wsBar.Range("G4").Select()
